# Applies the "Added filtering options for the Component Analysis" edit:
# shifts the quarter-ahead error columns (Q0..Q9, i.e. B:K) one column to the
# right for each data row (rows 2-24), dropping the old rightmost column (K)
# and inserting a newly computed Q0 value into column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = New-Object 'object[][]' 23

$values[0] = @(-0.8260246319213993,-1.981077477814098,-0.4414436339245075,-0.7110977538981412,-0.159566049999028,-0.3854928109118805,-0.2776770955153309,-0.3740767161796326,0.1293341692733759,0.3878561111968251)
$values[1] = @(-0.1405287498260583,-0.4101828697996921,0.1413488340994211,-0.08457792681343129,0.0232377885831183,-0.0731618320811834,0.4302490533718251,0.6887709952952742,-1.03521158968304,-0.3815723733203383)
$values[2] = @(0.7815531178611421,0.5556263569482897,0.6634420723448393,0.5670424516805376,1.070453337133546,1.328975279056995,-0.3950073059213186,0.2586319104413827,0.6221758453498138,-0.2785999995298156)
$values[3] = @(0.59085446968239,0.4944548490180882,0.9978657344710967,1.256387676394546,-0.467594908583768,0.1860443077789333,0.5495882426873645,-0.351187602192265,0.7191879547613879,0.142128040744926)
$values[4] = @(0.9421746901075616,1.200696632031011,-0.5232859529473031,0.1303532634153982,0.4938971983238294,-0.4068786465558001,0.6634969103978527,0.08643699638139091,-0.1191689877116253,0.8599759940094719)
$values[5] = @(-0.8773451131212686,-0.2237058967585673,0.1398380381498639,-0.7609378067297656,0.3094377502238873,-0.2676221637925746,-0.4732281478855908,0.5059168338355065,-0.2304295095018555,-0.2483993241560489)
$values[6] = @(0.04951982135029001,-0.8512560235293395,0.2191195334243133,-0.3579403805921485,-0.5635463646851647,0.4155986170359325,-0.3207477263014294,-0.3387175409556228,-0.315312675415286,-0.6804871111820268)
$values[7] = @(0.2874422675287736,-0.2896176464876882,-0.4952236305807045,0.4839213511403928,-0.2524249921969692,-0.2703948068511625,-0.2469899413108257,-0.6121643770775664,0.04365514009338567,-0.08721574084803801)
$values[8] = @(-0.5938852086106097,0.3852597731104875,-0.3510865702268745,-0.3690563848810678,-0.345651519340731,-0.7108259551074718,-0.0550064379365196,-0.1858773188779433,-0.2953951178020652,-1.125172904869618)
$values[9] = @(-0.2740769919852016,-0.292046806639395,-0.2686419410990581,-0.6338163768657989,0.02200314030515318,-0.1088677406362705,-0.2183855395603924,-1.048163326627946,-0.4676375374592979,-0.5510236388890779)
$values[10] = @(-0.09712588908289738,-0.4623003248496382,0.193519192321314,0.0626483113798903,-0.04686948754423159,-0.8766472746117848,-0.2961214854431372,-0.3795075868729171,-0.5974174644126757,0.9663603392234945)
$values[11] = @(0.1560540901775642,0.02518320923614054,-0.08433458968798135,-0.9141123767555346,-0.3335865875868869,-0.4169726890166669,-0.6348825665564255,0.9288952370797448,-0.4070884690023188,1.033525947527809)
$values[12] = @(-0.106237102096024,-0.9360148891635773,-0.3554890999949296,-0.4388752014247095,-0.6567850789644682,0.9069927246717022,-0.4289909814103615,1.011623435119767,-0.0888206791410559,-0.2990950773107258)
$values[13] = @(-0.2673832081489602,-0.3507693095787401,-0.5686791871184989,0.9950986165176716,-0.3408850895643921,1.099729326965736,-0.0007147872950865053,-0.2109891854647564,0.387555017918974,0.06188637714375961)
$values[14] = @(-0.4749660705519546,1.088811733084216,-0.2471719729978479,1.19344244353228,0.09299832927145768,-0.1172760688982122,0.4812681344855182,0.1555994937103038,0.8727792950859075,2.568651570723608)
$values[15] = @(1.097015279354651,-0.2389684267274125,1.201645989802716,0.1012018755418931,-0.1090725226277768,0.4894716807559536,0.1638030399807392,0.8809828413563429,2.576855116994044,9.421101911918621)
$values[16] = @(-0.2415579873788807,1.199056429151248,0.09861231489042488,-0.111662083279245,0.4868821201044854,0.161213479329271,0.8783932807048747,2.574265556342576,9.418512351267154,-8.132141731834157)
$values[17] = @(1.205589761734299,0.1051456474734768,-0.1051287506961931,0.4934154526875373,0.1677468119123229,0.8849266132879265,2.580798888925627,9.425045683850206,-8.125608399251105,-0.4162117995949584)
$values[18] = @(-0.006071047505593896,-0.2163454456752638,0.3821987577084666,0.05653011693325222,0.7737099183088558,2.469582193946557,9.313828988871135,-8.236825094230175,-0.5274284945740291,1.055484166312883)
$values[19] = @(-0.2580915896621678,0.3404526137215625,0.01478397294634815,0.7319637743219518,2.427836049959653,9.272082844884231,-8.278571238217079,-0.5691746385609331,1.013738022325978,-1.952548970023277)
$values[20] = @(0.3657676764542774,0.04009903567906303,0.7572788370546667,2.453151112692368,9.297397907616945,-8.253256175484365,-0.5438595758282182,1.039053085058693,-1.927233907290562,0.3367972473739005)
$values[21] = @(0.04413770072197692,0.7613175020975806,2.457189777735282,9.30143657265986,-8.249217510441451,-0.5398209107853043,1.043091750101607,-1.923195242247649,0.3408359124168144,-0.1103182425099242)
$values[22] = @(0.6424403654065582,2.338312641044259,9.182559435968837,-8.368094647132473,-0.6586980474763267,0.9242146134105849,-2.042072378938671,0.221958775725792,-0.2291953792009466,-0.3695594427149207)

for ($i = 0; $i -lt 23; $i++) {
    $row = $i + 2
    $rowValues = $values[$i]
    for ($j = 0; $j -lt 10; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
